$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new attendance record as the next row after the existing data
# (the sheet currently has data through row 71, so the new entry goes in row 72).
$row = 72

# Force columns to text before assigning, so a date-looking string like
# "2025-09-17" is stored as literal text (matching the rest of the sheet)
# instead of being auto-converted into a date serial number.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 3).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "Kartikey Gupta"
$ws.Cells.Item($row, 2).Value = "2025-09-17"
$ws.Cells.Item($row, 3).Value = "19:00:22"

# Reset the cell style back to Normal (default) so no explicit style index
# is left on the new cells, matching the plain unstyled data rows above.
$ws.Range("A72:C72").Style = "Normal"
